# Master_NCI_Risk_Table.xlsx update:
# 1. "HPV-positive NILM" -> "HPV-positive NILM x2" (sheet "3 - Colposcopy Results")
# 2. " HPV-negative/ASCUS/LSIL" -> " HPV-negative/ASC-US/LSIL" (sheet "4 - Post Colpo Surveillance")
# 3. Add a new "CIN 2 or 3" category and use it (instead of "CIN3") for the
#    "Biopsy Result Before Treatment" column on sheet "5 - Post Treatment Surveillance"

$wb = $excel.ActiveWorkbook

$wsColposcopy = $wb.Worksheets.Item("3 - Colposcopy Results")
$wsColposcopy.Cells.Replace("HPV-positive NILM", "HPV-positive NILM x2", 1)

$wsPostColpo = $wb.Worksheets.Item("4 - Post Colpo Surveillance")
$wsPostColpo.Cells.Replace(" HPV-negative/ASCUS/LSIL", " HPV-negative/ASC-US/LSIL", 1)

$wsPostTx = $wb.Worksheets.Item("5 - Post Treatment Surveillance")
for ($r = 2; $r -le 12; $r++) {
    $wsPostTx.Cells.Item($r, 1).Value = "CIN 2 or 3"
}
